$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells we touch are stored as Text so exact
# formatting (trailing zeros, leading zeros, etc.) is preserved.
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D23', 'D24', 'D26', 'D27', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D47', 'D49')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '247.15'

# Row 3
$ws.Range("D3").Value = '22.03'

# Row 4
$ws.Range("D4").Value = '5.464'

# Row 5
$ws.Range("D5").Value = '0.05777'

# Row 7
$ws.Range("D7").Value = '6.369'

# Row 8
$ws.Range("D8").Value = '0.8188'

# Row 9
$ws.Range("D9").Value = '0.9763'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1431'
$ws.Range("E10").Value = '9WazirXWRX'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07457'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03145'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.02998'
$ws.Range("E13").Value = '12BitrueCoinBTR'

# Row 14
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = '4.162'
$ws.Range("E14").Value = '13MCDexMCB'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09421'
$ws.Range("E15").Value = '14BitMartTokenBMX'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001587'
$ws.Range("E16").Value = '15BitForexTokenBF'

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04815'
$ws.Range("E17").Value = '16CoinExTokenCET'

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005849'
$ws.Range("E18").Value = '17OneONE'

# Row 19
$ws.Range("D19").Value = '0.006192'

# Row 20
$ws.Range("D20").Value = '0.004133'

# Row 21
$ws.Range("D21").Value = '0.0009946'

# Row 23
$ws.Range("D23").Value = '3.765'

# Row 24
$ws.Range("D24").Value = '2.221'

# Row 26
$ws.Range("D26").Value = '0.1260'

# Row 27
$ws.Range("D27").Value = '0.0003999'

# Row 40
$ws.Range("D40").Value = '0.03894'

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1077'
$ws.Range("E41").Value = '40BKEXTokenBKK'

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.002629'
$ws.Range("E42").Value = '41CEJICEJIBestin24h'

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '0.003062'
$ws.Range("E43").Value = '42KickTokenKICK'

# Row 44
$ws.Range("D44").Value = '0.006247'

# Row 45
$ws.Range("D45").Value = '0.00005595'

# Row 47
$ws.Range("D47").Value = '0.3800'

# Row 49
$ws.Range("D49").Value = '0.00002100'

